# Add p-values to correlation analysis
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells H1/I1 with the p-value labels; copy G1's formatting
# (bold font, borders, centered alignment) onto them so they match the rest
# of the header row.
$ws.Range("H1").Value2 = "Kendall's p-Value"
$ws.Range("I1").Value2 = "Spearman's p-Value"

$ws.Range("G1").Copy() | Out-Null
$ws.Range("H1:I1").PasteSpecial(-4122) | Out-Null

# Row 2 (TEMP placeholder row) - fill H2/I2 with TEMP same as F2/G2
$ws.Range("H2").Value2 = $ws.Range("G2").Value2
$ws.Range("I2").Value2 = $ws.Range("G2").Value2

# Row 3: the Spearman's Rho value moves from G3 to H3, new p-values fill G3 and I3
$ws.Range("H3").Value2 = $ws.Range("G3").Value2
$ws.Range("G3").Value2 = 0.00681274956405711
$ws.Range("I3").Value2 = 0.00611470670329826

# Row 4 (TEMP placeholder row) - fill H4/I4 with TEMP same as F4/G4
$ws.Range("H4").Value2 = $ws.Range("G4").Value2
$ws.Range("I4").Value2 = $ws.Range("G4").Value2
